$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Assign owner "Keye Li" and status "in progress" to the newly-added
# project rows (A99:A118).
for ($r = 99; $r -le 118; $r++) {
    $ws.Cells.Item($r, 2).Value = "Keye Li"
    $ws.Cells.Item($r, 3).Value = "in progress"
}

# Reflect the scroll position / active selection captured in the saved
# workbook view.
$ws.Range("E112").Select()
$excel.ActiveWindow.ScrollRow = 100
